$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the note text in D22 (more context added to the end of the existing note)
$ws.Range("D22").Value = "Finished mf using recosystem, not very good.. Looking to optimize. Workday banking setup and courses."

# Hours worked on 2024-05-24 corrected from 5 to 7
$ws.Range("B22").Value = 7

# That note now wraps across three lines, same as the other long notes in the sheet
$ws.Rows.Item(22).RowHeight = 45

# New entry for 2024-05-26
$ws.Range("A23").Value = 45438
$ws.Range("A23").NumberFormat = $ws.Range("A22").NumberFormat
$ws.Range("B23").Value = 3
$ws.Range("D23").Value = "Workday courses and minor testing"
$ws.Range("F23").Value = "Changing course, no MF, at least not reguarly. Back to the drawing board."

# Move the view/selection down to the newly added row, as the author left it
$ws.Range("B24").Select()
$excel.ActiveWindow.ScrollRow = 14
